$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update I20 (days) and K20 (Friday date) on row 20
$ws.Range("I20").Value = 4
$ws.Range("K20").Value = 2

# Restore the normal (non-greyed) border style on G20 to match D20:F20
$ws.Range("F20").Copy()
$ws.Range("G20").PasteSpecial(-4122)  # xlPasteFormats

# Update the K column "Friday date" values for rows 21-24
$ws.Range("K21").Value = 9
$ws.Range("K22").Value = 16
$ws.Range("K23").Value = 23
$ws.Range("K24").Value = 30

# Update the active cell selection to L20
$ws.Range("L20").Select()
